$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion rates text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 1.79 = 6472.27 pesos", "1000 Bs = 1.85 = 6699.86 pesos")
$text = $text.Replace("6472.27 pesos = 1.78 = 935.68 Bs", "6699.86 pesos = 1.85 = 971.51 Bs")
$cell.Value = $text

# --- Sheet "tasas": update the updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 540.3099999999999
$ws2.Range("O10").Value = 3620
$ws2.Range("N12").Value = 3630
$ws2.Range("O12").Value = 526.369
